# Update the "想去人数" (interested count) figures in column F for the
# events whose numbers changed between scrapes, on both the "展览"
# sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 283
    7  = 1345
    8  = 480
    10 = 165
    11 = 118
    12 = 168
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
